# Update "想去人数" (interest count) figures across sheets.
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 164
$ws1.Range("F3").Value = 479
$ws1.Range("F4").Value = 15
$ws1.Range("F5").Value = 19
$ws1.Range("F8").Value = 17
$ws1.Range("F9").Value = 573

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 86
$ws2.Range("F3").Value = 37

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 164
$ws4.Range("F3").Value = 86
$ws4.Range("F4").Value = 479
$ws4.Range("F5").Value = 15
$ws4.Range("F6").Value = 19
$ws4.Range("F9").Value = 17
$ws4.Range("F10").Value = 573
$ws4.Range("F11").Value = 37
